{"js": "// Append additional commentary to the \"\u010clani projekta\" paragraph that\n// currently ends with \"..., Jan Vrta\u010dnik ter Alja\u017e Mar\u0161. \".\n// The new text is appended as additional runs inside the SAME paragraph\n// (no new paragraph is created), matching the source OOXML diff.\n\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items/text\");\nawait context.sync();\n\nconst marker = \"Jan Vrta\u010dnik ter Alja\u017e Mar\u0161.\";\nlet target = null;\nfor (let i = 0; i < paragraphs.items.length; i++) {\n  const p = paragraphs.items[i];\n  if (p.text && p.text.indexOf(marker) !== -1) {\n    target = p;\n    break;\n  }\n}\n\nif (!target) {\n  throw new Error(\"Could not locate the 'Jan Vrta\u010dnik ter Alja\u017e Mar\u0161.' paragraph.\");\n}\n\nconst addition =\n  \"Pri delu smo si sicer pomagali med sabo, vendar je za izdelavo aplikacije za telefon, ki je uporabljena za zbiranje podatkov, bil glavni Jan Vrta\u010dnik. Jakob Polegek in Alja\u017e Mar\u0161 pa sta bila zadol\u017eena za izdelavo spletne strani, ki ima glavno vlogo prikaz podatkov. Jakob Polegek je prav tako delal na aplikaciji, ki pridobi iz videa uporabnika, njegov sr\u010dni utrip. Alja\u017e Mar\u0161 pa je bil zadol\u017een za ustvarjanje spletnega stre\u017enika, ki bo kodo, za pridobivanje povpre\u010dnega sr\u010dnega utripa lahko poganjala. Pri izbiri podatkovne baze pa so bili prisotni vsi \u010dlani skupine.\";\n\ntarget.insertText(addition, Word.InsertLocation.end);\nawait context.sync();\n", "ps1": "# Append additional commentary to the \"\u010clani projekta\" paragraph that\n# currently ends with \"..., Jan Vrta\u010dnik ter Alja\u017e Mar\u0161. \".\n# The new text is appended as additional text inside the SAME paragraph\n# (no new paragraph is created), matching the source OOXML diff.\n\n$d = $word.ActiveDocument\n\n$addition = \"Pri delu smo si sicer pomagali med sabo, vendar je za izdelavo aplikacije za telefon, ki je uporabljena za zbiranje podatkov, bil glavni Jan Vrta\u010dnik. Jakob Polegek in Alja\u017e Mar\u0161 pa sta bila zadol\u017eena za izdelavo spletne strani, ki ima glavno vlogo prikaz podatkov. Jakob Polegek je prav tako delal na aplikaciji, ki pridobi iz videa uporabnika, njegov sr\u010dni utrip. Alja\u017e Mar\u0161 pa je bil zadol\u017een za ustvarjanje spletnega stre\u017enika, ki bo kodo, za pridobivanje povpre\u010dnega sr\u010dnega utripa lahko poganjala. Pri izbiri podatkovne baze pa so bili prisotni vsi \u010dlani skupine.\"\n\n$found = $false\nforeach ($p in $d.Paragraphs) {\n    $r = $p.Range\n    if ($r.Text -like \"*Jan Vrta\u010dnik ter Alja\u017e Mar\u0161.*\") {\n        $insertRange = $r.Duplicate\n        $insertRange.Collapse(0)\n        [void]$insertRange.MoveEnd(1, -1)\n        $insertRange.InsertAfter($addition)\n        $found = $true\n        break\n    }\n}\n\nif (-not $found) {\n    throw \"Could not locate the 'Jan Vrta\u010dnik ter Alja\u017e Mar\u0161.' paragraph.\"\n}\n"}
